# Applies the "fixed PER bug" edit:
#   - Column A (index) values stay the same for each row.
#   - Column B (team) and Column C (value) are updated row-by-row to reflect
#     the corrected / re-derived dataset.
#
# New per-row values (row -> Team, Value):
$updates = @(
    @{Row=2;  Team="POR"; Value=12.24285714285714},
    @{Row=3;  Team="CLE"; Value=13.49285714285714},
    @{Row=4;  Team="DAL"; Value=14.82857142857143},
    @{Row=5;  Team="MIA"; Value=14.725},
    @{Row=6;  Team="OKC"; Value=13.34285714285714},
    @{Row=7;  Team="ATL"; Value=14.00769230769231},
    @{Row=8;  Team="WAS"; Value=12.66428571428572},
    @{Row=9;  Team="MIL"; Value=12.26666666666667},
    @{Row=10; Team="LAC"; Value=14.12857142857143},
    @{Row=11; Team="SAS"; Value=16.85714285714285},
    @{Row=12; Team="DET"; Value=12.75},
    @{Row=13; Team="ORL"; Value=12.13076923076923},
    @{Row=14; Team="UTA"; Value=12.23333333333333},
    @{Row=15; Team="MEM"; Value=13.84285714285715},
    @{Row=16; Team="HOU"; Value=13.72},
    @{Row=17; Team="NOP"; Value=12.675},
    @{Row=18; Team="DEN"; Value=13.33846153846154},
    @{Row=19; Team="LAL"; Value=12.11333333333334},
    @{Row=20; Team="GSW"; Value=16.15714285714285},
    @{Row=21; Team="IND"; Value=14.86428571428572},
    @{Row=22; Team="CHO"; Value=14.98571428571429},
    @{Row=23; Team="CHI"; Value=13.54285714285715},
    @{Row=24; Team="PHI"; Value=12.69285714285714},
    @{Row=25; Team="BOS"; Value=13.70666666666667},
    @{Row=26; Team="BRK"; Value=13.82},
    @{Row=27; Team="TOR"; Value=13.56},
    @{Row=28; Team="SAC"; Value=12.77333333333334},
    @{Row=29; Team="PHO"; Value=13.67857142857143},
    @{Row=30; Team="NYK"; Value=14.03333333333333},
    @{Row=31; Team="MIN"; Value=12.12142857142857}
)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 2).Value = $u.Team
    $ws.Cells.Item($u.Row, 3).Value = $u.Value
}
